$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Existing rows 2-60 (years 1950-2008): update the Data (column E) values in place
# New rows 61-68 (years 2009-2016): append new rows for Libya GDP per Capita

$countryCode = 434
$countryName = "Libya"
$indicator = "GDP per Capita"

$yearData = @{
    1950 = "620"
    1951 = "668"
    1952 = "697"
    1953 = "685"
    1954 = "662"
    1955 = "807"
    1956 = "950"
    1957 = "958"
    1958 = "1020"
    1959 = "1063"
    1960 = "1323"
    1961 = "1400"
    1962 = "1788"
    1963 = "2322"
    1964 = "3156"
    1965 = "3889"
    1966 = "4417"
    1967 = "4733"
    1968 = "6071"
    1969 = "6539"
    1970 = "6591"
    1971 = "6008"
    1972 = "6237"
    1973 = "6005"
    1974 = "6923"
    1975 = "6872"
    1976 = "8142"
    1977 = "8682"
    1978 = "8730"
    1979 = "9086"
    1980 = "8741"
    1981 = "6760"
    1982 = "6661"
    1983 = "6237"
    1984 = "5695"
    1985 = "6092"
    1986 = "5370"
    1987 = "4465"
    1988 = "4670"
    1989 = "4873"
    1990 = "4921"
    1991 = "5977.33692478306"
    1992 = "6113.35350838343"
    1993 = "6190.68138314581"
    1994 = "6645.07657054471"
    1995 = "6898.95237639788"
    1996 = "7543.93543073536"
    1997 = "8422.2146757838"
    1998 = "8545.48067270692"
    1999 = "9064.06731938906"
    2000 = "9770.98043889931"
    2001 = "10354.4825766246"
    2002 = "10787.6747325271"
    2003 = "12870.5745261769"
    2004 = "14193.2590078302"
    2005 = "16538.5566293337"
    2006 = "18657.0616505443"
    2007 = "20731.7177828003"
    2008 = "25140.5988288157"
    2009 = "26403.7880429989"
    2010 = "29157.1420833813"
    2011 = "12893"
    2012 = "29766"
    2013 = "13326"
    2014 = "9736"
    2015 = "8514"
    2016 = "8096"
}

$rng = $ws.Range("E2:E68")
$rng.NumberFormat = "@"

for ($row = 2; $row -le 68; $row++) {
    $year = 1948 + $row
    $ws.Cells.Item($row, 1).Value = $countryCode
    $ws.Cells.Item($row, 2).Value = $countryName
    $ws.Cells.Item($row, 3).Value = $indicator
    $ws.Cells.Item($row, 4).Value = $year
    $ws.Cells.Item($row, 5).Value = $yearData[$year]
}

$rng.ClearFormats()
